# Auto-generated Excel COM-interop script
# Applies the cryptos.xlsx update described by the commit diff:
#  - Rows 2-33: refresh Price (D) and Volume(1h) (E) figures only
#  - A new "Frax" entry is inserted as row 34, so rows 34-50 shift down
#    to rows 35-51 with refreshed Price/Volume values
#  - The former row 51 ("Aave") drops off the bottom of the A1:E51 table
#
# Price values in column D are plain text (e.g. "0.9990", "26.450.78")
# in the source workbook, so a helper is used to force text storage and
# avoid Excel auto-converting them into numbers (which would also drop
# significant trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

# --- Update existing rows 2-33 (Price / Volume columns only) ---
Set-TextValue $ws.Range("D2") "26.450.78"
$ws.Range("E2").Value = "  -0.33%  "

Set-TextValue $ws.Range("D3") "1.725.66"
$ws.Range("E3").Value = "  -0.13%  "

Set-TextValue $ws.Range("D4") "0.9990"
$ws.Range("E4").Value = "  -0.09%  "

Set-TextValue $ws.Range("D5") "243.32"
$ws.Range("E5").Value = "  -0.81%  "

$ws.Range("E6").Value = "  -0.07%  "

Set-TextValue $ws.Range("D7") "0.4921"
$ws.Range("E7").Value = "  +2.46%  "

Set-TextValue $ws.Range("D8") "0.2621"
$ws.Range("E8").Value = "  -1.53%  "

Set-TextValue $ws.Range("D9") "0.06207"
$ws.Range("E9").Value = "  +0.54%  "

Set-TextValue $ws.Range("D10") "1.725.72"
$ws.Range("E10").Value = "  -0.09%  "

Set-TextValue $ws.Range("D11") "0.06993"
$ws.Range("E11").Value = "  -2.55%  "

Set-TextValue $ws.Range("D12") "15.44"
$ws.Range("E12").Value = "  -0.63%  "

Set-TextValue $ws.Range("D13") "4.568"
$ws.Range("E13").Value = "  +1.05%  "

Set-TextValue $ws.Range("D14") "0.6010"
$ws.Range("E14").Value = "  -0.98%  "

Set-TextValue $ws.Range("D15") "77.38"
$ws.Range("E15").Value = "  +0.38%  "

Set-TextValue $ws.Range("D16") "0.9994"
$ws.Range("E16").Value = "  -0.11%  "

Set-TextValue $ws.Range("D17") "26.446.32"
$ws.Range("E17").Value = "  -0.35%  "

Set-TextValue $ws.Range("D18") "0.9992"
$ws.Range("E18").Value = "  -0.15%  "

Set-TextValue $ws.Range("D19") "0.000007190"
$ws.Range("E19").Value = "  +3.38%  "

Set-TextValue $ws.Range("D20") "11.36"
$ws.Range("E20").Value = "  -1.45%  "

Set-TextValue $ws.Range("D21") "1.948.88"
$ws.Range("E21").Value = "  -0.24%  "

Set-TextValue $ws.Range("D22") "4.489"
$ws.Range("E22").Value = "  -0.65%  "

Set-TextValue $ws.Range("D23") "8.601"
$ws.Range("E23").Value = "  -2.05%  "

Set-TextValue $ws.Range("D24") "5.166"
$ws.Range("E24").Value = "  -1.44%  "

Set-TextValue $ws.Range("D25") "138.25"
$ws.Range("E25").Value = "  +1.03%  "

Set-TextValue $ws.Range("D26") "15.28"
$ws.Range("E26").Value = "  -0.58%  "

Set-TextValue $ws.Range("D27") "1.397"
$ws.Range("E27").Value = "  -0.56%  "

Set-TextValue $ws.Range("D28") "107.16"
$ws.Range("E28").Value = "  -0.10%  "

Set-TextValue $ws.Range("D29") "1.721"
$ws.Range("E29").Value = "  -3.11%  "

Set-TextValue $ws.Range("D30") "3.960"
$ws.Range("E30").Value = "  -0.21%  "

Set-TextValue $ws.Range("D31") "0.07987"
$ws.Range("E31").Value = "  -0.19%  "

Set-TextValue $ws.Range("D32") "3.681"
$ws.Range("E32").Value = "  -0.03%  "

Set-TextValue $ws.Range("D33") "0.04524"
$ws.Range("E33").Value = "  +0.47%  "

# --- Rows 34-51: data shifted down by one due to new "Frax" row; update Coin/Link/Price/Volume ---
$ws.Range("B34").Value = "Frax"
$ws.Range("C34").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D34") "0.9987"
$ws.Range("E34").Value = "  -0.14%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D35") "2.601"
$ws.Range("E35").Value = "  -0.61%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D36") "0.9995"
$ws.Range("E36").Value = "  -0.22%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D37") "0.6276"
$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D38") "0.9330"
$ws.Range("E38").Value = "  +2.66%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D39") "1.964"
$ws.Range("E39").Value = "  -4.40%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D40") "2.390"
$ws.Range("E40").Value = "  -0.30%  "

$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D41") "0.9995"
$ws.Range("E41").Value = "  -0.42%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D42") "0.01490"
$ws.Range("E42").Value = "  -0.67%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D43") "99.60"
$ws.Range("E43").Value = "  -2.78%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D44") "5.359"
$ws.Range("E44").Value = "  -2.59%  "

$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D45") "0.3856"
$ws.Range("E45").Value = "  -0.64%  "

$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D46") "6.791"
$ws.Range("E46").Value = "  -3.48%  "

$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D47") "0.1168"
$ws.Range("E47").Value = "  -0.87%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D48") "0.05370"
$ws.Range("E48").Value = "  -0.14%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D49") "7.773"
$ws.Range("E49").Value = "  -0.45%  "

$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue $ws.Range("D50") "30.22"
$ws.Range("E50").Value = "  -1.54%  "

$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D51") "1.231"
$ws.Range("E51").Value = "  -1.11%  "

